# Books.xlsx edit: split the "Количество томов в год" (volumes-in-year count)
# column into two explicit columns - first/last volume number in the year -
# and move the old "Том (старая маркировка)" column to the end (new column H).
#
# Old layout (A..G): Автор | Название | Том | Количество томов | Год издания |
#                     Том (старая маркировка) | Количество томов в год
# New layout (A..H): Автор | Название | Том | Количество томов | Год издания |
#                     Номер первого тома в году | Номер последнего тома в году |
#                     Том (старая маркировка)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Move the old column-F data (rows that actually have a value: the header
#     and the Brockhaus/Efron encyclopedia rows 157-163) into the new column H
#     before column F gets overwritten with the new "first volume" numbers.
$oldF = @{}
for ($r = 1; $r -le 179; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $oldF[$r] = $val
    }
}
foreach ($r in $oldF.Keys) {
    $ws.Cells.Item($r, 8).Value() = $oldF[$r]
}

# --- Header row (row 1) ---
$ws.Range("F1").Value() = "Номер первого тома в году"
$ws.Range("G1").Value() = "Номер последнего тома в году"
# H1 already received "Том (старая маркировка)" via the copy loop above.

# --- Data rows 157-163 (Энциклопедия Брокгауза и Ефрона) ---
# For each row, replace the old "count of volumes in year" (column G) with an
# explicit first/last volume-number pair in columns F/G.
$firstVol = @{157=4; 158=16; 159=45; 160=56; 161=62; 162=68; 163=73}
$lastVol  = @{157=9; 158=21; 159=50; 160=61; 161=67; 162=72; 163=78}

foreach ($r in $firstVol.Keys) {
    $ws.Cells.Item($r, 6).Value() = $firstVol[$r]
    $ws.Cells.Item($r, 7).Value() = $lastVol[$r]
}

# --- Column widths (match bestFit widths produced by the new, wider header
#     text in F/G and the relocated column H) ---
$ws.Range("F1:H163").EntireColumn.AutoFit() | Out-Null

# --- Selection / view state ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 137
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G162").Select() | Out-Null
